$d = $word.ActiveDocument

# 1. Update the wording about the traceback step / complexity discussion.
#    Old:  "...iterates over every city. If we had instead..."
#    New:  "...iterates over every city, however adding the traceback step
#           yields the final result. If we had instead..."
$d.Content.Find.Execute(
    "every city. If we had instead",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "every city, however adding the traceback step yields the final result. If we had instead",
    2
)

# 2. Move the "_GoBack" last-edit bookmark so it sits right after the newly
#    typed text (where Word would leave it following a live edit), i.e.
#    immediately before "yields the final result.".
$bm = $d.Content
$bm.Find.Execute(" step ")
$bm.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bm)
